$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D; this shifts existing D:K data to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy the number formatting from column F (the old column D, now shifted) into
# the two newly inserted columns D:E so the new cells inherit the correct
# date / number styles for every row.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)

# Populate the two new quarters of data (columns D and E).
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 41400
$ws.Range("E8").Value = 39000
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = -300
$ws.Range("E15").Value = -300
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("D17").Value = 13000
$ws.Range("E17").Value = 10300
$ws.Range("D18").Value = 28400
$ws.Range("E18").Value = 28700
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("D20").Value = -18700
$ws.Range("E20").Value = -15600
$ws.Range("D21").Value = 11300
$ws.Range("E21").Value = 14800
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 9700
$ws.Range("E23").Value = 13100
$ws.Range("D24").Value = 2200
$ws.Range("E24").Value = 2600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 7500
$ws.Range("E26").Value = 10600
$ws.Range("D27").Value = 7100
$ws.Range("E27").Value = 10200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 18700
$ws.Range("E32").Value = 15600
$ws.Range("D33").Value = 7100
$ws.Range("E33").Value = 10200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 7100
$ws.Range("E35").Value = 10200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D39").Value = ""
$ws.Range("E39").Value = ""
$ws.Range("D40").Value = ""
$ws.Range("E40").Value = ""
$ws.Range("D41").Value = 102800
$ws.Range("E41").Value = 117300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 42800
$ws.Range("E48").Value = 43300
$ws.Range("D49").Value = 76200
$ws.Range("E49").Value = 78900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 4311700
$ws.Range("E54").Value = 4258400
$ws.Range("D55").Value = ""
$ws.Range("E55").Value = ""
$ws.Range("D56").Value = ""
$ws.Range("E56").Value = ""
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 39200
$ws.Range("E61").Value = 39200
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 3915400
$ws.Range("E66").Value = 3866200
$ws.Range("D67").Value = ""
$ws.Range("E67").Value = ""
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 17300
$ws.Range("E70").Value = 17300
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 279900
$ws.Range("E72").Value = 276600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 379000
$ws.Range("E76").Value = 374800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 7100
$ws.Range("E81").Value = 10200
$ws.Range("D82").Value = ""
$ws.Range("E82").Value = ""
$ws.Range("D83").Value = 1600
$ws.Range("E83").Value = 1600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 14300
$ws.Range("E89").Value = 20500
$ws.Range("D90").Value = ""
$ws.Range("E90").Value = ""
$ws.Range("D91").Value = -800
$ws.Range("E91").Value = -900
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -67200
$ws.Range("E94").Value = -47100
$ws.Range("D95").Value = ""
$ws.Range("E95").Value = ""
$ws.Range("D96").Value = -4200
$ws.Range("E96").Value = -4200
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 38300
$ws.Range("E100").Value = 54800
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -14600
$ws.Range("E102").Value = 28200
